$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.4800474766800295
$ws.Range("J2").Value = 0.4800474766800294
$ws.Range("M2").Value = 0.08849299999999999
$ws.Range("R2").Value = 0.5597340652469999
$ws.Range("S2").Value = 0.4345540115031014
$ws.Range("T2").Value = 0.4345540115031013

$ws.Range("I3").Value = 0.4800474766800295
$ws.Range("J3").Value = 0.4800474766800294
$ws.Range("O3").Value = 0.09476867890558938
$ws.Range("P3").Value = 0.09476867890558936
$ws.Range("S3").Value = 0.04549346517692811
$ws.Range("T3").Value = 0.0454934651769281

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.341519
$ws.Range("H4").Value = 1.024557
$ws.Range("I4").Value = 0.233275296666637
$ws.Range("J4").Value = 0.233275296666637
$ws.Range("M4").Value = 0.08849299999999999
$ws.Range("Q4").Value = 0.03022204086699999
$ws.Range("R4").Value = 0.271998367803
$ws.Range("S4").Value = 0.2111681049802304
$ws.Range("T4").Value = 0.2111681049802304

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.341519
$ws.Range("H5").Value = 1.024557
$ws.Range("I5").Value = 0.233275296666637
$ws.Range("J5").Value = 0.233275296666637
$ws.Range("O5").Value = 0.09476867890558938
$ws.Range("P5").Value = 0.09476867890558936
$ws.Range("Q5").Value = 0.003163945855666667
$ws.Range("R5").Value = 0.028475512701
$ws.Range("S5").Value = 0.02210719168640663
$ws.Range("T5").Value = 0.02210719168640663

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4197003333333333
$ws.Range("H6").Value = 1.259101
$ws.Range("I6").Value = 0.2866772266533335
$ws.Range("J6").Value = 0.2866772266533335
$ws.Range("M6").Value = 0.08849299999999999
$ws.Range("Q6").Value = 0.03714054159766667
$ws.Range("R6").Value = 0.334264874379
$ws.Range("S6").Value = 0.2595092046110788
$ws.Range("T6").Value = 0.2595092046110788

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4197003333333333
$ws.Range("H7").Value = 1.259101
$ws.Range("I7").Value = 0.2866772266533335
$ws.Range("J7").Value = 0.2866772266533335
$ws.Range("O7").Value = 0.09476867890558938
$ws.Range("P7").Value = 0.09476867890558936
$ws.Range("Q7").Value = 0.003888243788111112
$ws.Range("R7").Value = 0.034994194093
$ws.Range("S7").Value = 0.02716802204225463
$ws.Range("T7").Value = 0.02716802204225462
